$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Shape "Zylinder 3": "<Tbd>" -> "Postgres" + " in Memory" ---
$shape1 = $s.Shapes.Item(1)
$tr1 = $shape1.TextFrame.TextRange
$para2 = $tr1.Paragraphs(2, 1)

# Run 1 is "<" with plain rPr -> drop it, its text moves into the
# ("Tbd", err="1") run so that run keeps its spell-check flag on "Postgres".
$run1 = $para2.Runs(1, 1)
$run1.Text = ""

$para2 = $tr1.Paragraphs(2, 1)
$run2 = $para2.Runs(1, 1)
$run2.Text = "Postgres"

$para2 = $tr1.Paragraphs(2, 1)
$run3 = $para2.Runs(2, 1)
$run3.Text = " in Memory"

# --- Shape "Pfeil: nach links und rechts 12": "<tbd>, Spring Data?" -> "Spring Data & JPA" ---
$shape2 = $s.Shapes.Item(8)
$tr2 = $shape2.TextFrame.TextRange
$paraA = $tr2.Paragraphs(1, 1)

$runA1 = $paraA.Runs(1, 1)
$runA1.Text = "Spring Data & JPA"

$paraA = $tr2.Paragraphs(1, 1)
$runA2 = $paraA.Runs(2, 1)
$runA2.Text = ""

$paraA = $tr2.Paragraphs(1, 1)
$runA3 = $paraA.Runs(2, 1)
$runA3.Text = ""
